# Replace the metabolite row-name labels ("Met1", "Met2", "Met3") with
# plain string row indices ("1", "2", "3") on every worksheet that uses
# them, making sure the values remain text (not numbers).

$wb = $excel.ActiveWorkbook

$oldNames = @("Met1", "Met2", "Met3")
$newNames = @("1", "2", "3")

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $current = $cell.Value2
            for ($i = 0; $i -lt $oldNames.Length; $i++) {
                if ($current -eq $oldNames[$i]) {
                    # Keep the cell formatted as text and write the new
                    # value as a string so it is stored as a shared string
                    # (not re-interpreted as a number) in the saved file.
                    $cell.NumberFormat = "@"
                    $cell.Value = $newNames[$i]
                }
            }
        }
    }
}
